# Ohio Vaccination URLs workbook update
# - Rewrites the B24:B26 URLs to point at raw.githubusercontent.com instead
#   of github.com/.../blob/...  (one hash also changed upstream: row 25).
# - Appends 12 new date/URL rows (27-38).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 24-26 (URL text only; keep date + layout) -------
$ws.Range("B24").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/f572912a4f5ab611928e37836fef3c7fff775ab6/incidents.json"
$ws.Range("B25").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/31c91307b2438eca531a836aa3e6bd4e3efd129c/incidents.json"
$ws.Range("B26").Value = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/4ea90d437fe37aade429c0d3550050c532d5f729/incidents.json"

# --- New rows 27-38: dates + URLs ------------------------------------------
$dates = @(44279, 44280, 44281, 44282, 44283, 44284, 44285, 44286, 44287, 44288, 44289, 44290)
$urls = @(
  "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/266081e25c58d9110b33c9fa0bdeafb546dcadbb/incidents.json",
  "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/8e6ebfe0cf10d5baa7bd29257e51b905a760819a/incidents.json",
  "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/01781c4f56c7a408cd6753e7bc082404c5d960ee/incidents.json",
  "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/58815bc9990352b81277064b97ddc32f60eca836/incidents.json",
  "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/eb0b1643db37830f1f216ca41abaffa9e7b77089/states.json",
  "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/299aabaccc94b4728cf350df468ba65fa5fbbb87/states.json",
  "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/776a1964a807016e46c5cde11cd707716bc4a5aa/states.json",
  "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/bdbedcab55d3f575ad68553f347891229264c778/states.json",
  "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/f3fc91f1885a0d5fa415bcf4d6bc7a1c9b9b2f2f/states.json",
  "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/13bf1552463c9afa51b12cab8aad6c3a36d416f0/states.json",
  "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/f89e0c8166b25645f8e92ebc148f9cb9db119554/states.json",
  "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/0fc524b44b024cefc5d3ce9407cbe2027ae93e79/states.json"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
  $row = 27 + $i
  # Carry the existing row's formatting (date number format on col A,
  # plain text on col B) down onto the freshly-appended row.
  $ws.Range("A26:B26").Copy()
  $ws.Range("A" + $row + ":B" + $row).PasteSpecial(-4122)

  $ws.Range("A" + $row).Value = $dates[$i]
  $ws.Range("B" + $row).Value = $urls[$i]
}

$excel.CutCopyMode = 0

# Match the saved view-state from the source session (scrolled back to the
# top, cursor left on F7).
$ws.Range("F7").Select() | Out-Null
